$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2375.5715
$ws.Range("I2").Value = 328.8
$ws.Range("J2").Value = 7492.5
$ws.Range("K2").Value = 328.8
$ws.Range("L2").Value = 7492.5
$ws.Range("M2").Value = -215.8
$ws.Range("N2").Value = -7718.5

$ws.Range("H40").Value = 1948.4419
$ws.Range("I40").Value = 1765.3125
$ws.Range("J40").Value = 2481.182
$ws.Range("K40").Value = 1765.3125
$ws.Range("L40").Value = 2481.182
$ws.Range("M40").Value = -1590.3125
$ws.Range("N40").Value = -2831.182

$ws.Range("H51").Value = 5000
$ws.Range("I51").Value = 5000
$ws.Range("K51").Value = 5000
$ws.Range("M51").Value = -4516

$ws.Range("H88").Value = 1541.6923
$ws.Range("I88").Value = 2084.75
$ws.Range("J88").Value = 1300.3334
$ws.Range("K88").Value = 2084.75
$ws.Range("L88").Value = 1300.3334
$ws.Range("M88").Value = -1678.75
$ws.Range("N88").Value = -2112.3334

$ws.Range("H91").Value = 1541.6923
$ws.Range("I91").Value = 2084.75
$ws.Range("J91").Value = 1300.3334
$ws.Range("K91").Value = 2084.75
$ws.Range("L91").Value = 1300.3334
$ws.Range("M91").Value = -680.75
$ws.Range("N91").Value = -4108.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1364.8182
$ws.Range("I86").Value = 1171.625
$ws.Range("J86").Value = 1880
$ws.Range("K86").Value = 1171.625
$ws.Range("L86").Value = 1880
$ws.Range("M86").Value = -48.625
$ws.Range("N86").Value = -4126

$ws.Range("H89").Value = 1364.8182
$ws.Range("I89").Value = 1171.625
$ws.Range("J89").Value = 1880
$ws.Range("K89").Value = 5858.125
$ws.Range("L89").Value = 9400
$ws.Range("M89").Value = -242.125
$ws.Range("N89").Value = -20632

$ws.Range("H105").Value = 1554.8572
$ws.Range("I105").Value = 1414
$ws.Range("K105").Value = 1414
$ws.Range("M105").Value = 333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2277.0625
$ws.Range("I31").Value = 1940.7
$ws.Range("K31").Value = 1940.7
$ws.Range("M31").Value = -1645.7

$ws.Range("H34").Value = 2277.0625
$ws.Range("I34").Value = 1940.7
$ws.Range("K34").Value = 1940.7
$ws.Range("M34").Value = -1738.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 540.5
$ws.Range("J2").Value = 540.5
$ws.Range("L2").Value = 3243
$ws.Range("N2").Value = -3469

$ws.Range("H113").Value = 1548.1666
$ws.Range("J113").Value = 1498.5
$ws.Range("L113").Value = 4495.5
$ws.Range("N113").Value = -8835.5

$ws.Range("H121").Value = 18343.143
$ws.Range("I121").Value = 27567
$ws.Range("J121").Value = 6044.6665
$ws.Range("K121").Value = 82701
$ws.Range("L121").Value = 18133.9995
$ws.Range("M121").Value = -81391
$ws.Range("N121").Value = -20753.9995

$ws.Range("H131").Value = 627037.8
$ws.Range("I131").Value = 1638.75
$ws.Range("K131").Value = 4916.25
$ws.Range("M131").Value = 123.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 110
$ws.Range("J2").Value = 133.8
$ws.Range("L2").Value = 133.8
$ws.Range("N2").Value = -359.8

$ws.Range("H11").Value = 2456000
$ws.Range("I11").Value = 3062500
$ws.Range("K11").Value = 3062500
$ws.Range("M11").Value = -3062361

$ws.Range("H12").Value = 12500
$ws.Range("J12").Value = 12500
$ws.Range("L12").Value = 12500
$ws.Range("N12").Value = -12780

$ws.Range("H113").Value = 998.7778
$ws.Range("I113").Value = 998.7778
$ws.Range("K113").Value = 998.7778
$ws.Range("M113").Value = 1171.2222

$ws.Range("H122").Value = 3223.4666
$ws.Range("I122").Value = 3334.7693
$ws.Range("K122").Value = 10004.3079
$ws.Range("M122").Value = -7554.3079

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8974.3125
$ws.Range("I7").Value = 8359.6
$ws.Range("K7").Value = 8359.6
$ws.Range("M7").Value = -8247.6

$ws.Range("H20").Value = 8124.5
$ws.Range("J20").Value = 8124.5
$ws.Range("L20").Value = 8124.5
$ws.Range("N20").Value = -8576.5

$ws.Range("H26").Value = 10000
$ws.Range("I26").Value = 10000
$ws.Range("K26").Value = 10000
$ws.Range("M26").Value = -9705

$ws.Range("H29").Value = 20000
$ws.Range("I29").Value = 20000
$ws.Range("K29").Value = 20000
$ws.Range("M29").Value = -19705

$ws.Range("H40").Value = 4388.5557
$ws.Range("I40").Value = 3582.8333
$ws.Range("K40").Value = 3582.8333
$ws.Range("M40").Value = -3446.8333

$ws.Range("H55").Value = 851.6
$ws.Range("I55").Value = 649.7857
$ws.Range("J55").Value = 1108.4546
$ws.Range("K55").Value = 649.7857
$ws.Range("L55").Value = 1108.4546
$ws.Range("M55").Value = -476.7857
$ws.Range("N55").Value = -1454.4546

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H107").Value = 15000
$ws.Range("I107").Value = 15000
$ws.Range("K107").Value = 15000
$ws.Range("M107").Value = -13080

$ws.Range("H126").Value = 8974.3125
$ws.Range("I126").Value = 8359.6
$ws.Range("K126").Value = 25078.8
$ws.Range("M126").Value = -22608.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 26400
$ws.Range("J62").Value = 16499.5
$ws.Range("L62").Value = 16499.5
$ws.Range("N62").Value = -17747.5

$ws.Range("H65").Value = 26400
$ws.Range("J65").Value = 16499.5
$ws.Range("L65").Value = 82497.5
$ws.Range("N65").Value = -88737.5

$ws.Range("H74").Value = 18386.334
$ws.Range("I74").Value = 18377.5
$ws.Range("J74").Value = 18388.857
$ws.Range("K74").Value = 18377.5
$ws.Range("L74").Value = 18388.857
$ws.Range("M74").Value = -17441.5
$ws.Range("N74").Value = -20260.857

$ws.Range("H77").Value = 18386.334
$ws.Range("I77").Value = 18377.5
$ws.Range("J77").Value = 18388.857
$ws.Range("K77").Value = 55132.5
$ws.Range("L77").Value = 55166.571
$ws.Range("M77").Value = -50452.5
$ws.Range("N77").Value = -64526.571

$ws.Range("H122").Value = 1285.5714
$ws.Range("I122").Value = 880
$ws.Range("J122").Value = 2299.5
$ws.Range("K122").Value = 2640
$ws.Range("L122").Value = 6898.5
$ws.Range("M122").Value = -190
$ws.Range("N122").Value = -11798.5

$ws.Range("H126").Value = 2899.4285
$ws.Range("I126").Value = 2099.5
$ws.Range("K126").Value = 6298.5
$ws.Range("M126").Value = -3828.5

$ws.Range("H132").Value = 2290.4546
$ws.Range("I132").Value = 2319.5
$ws.Range("K132").Value = 6958.5
$ws.Range("M132").Value = -4428.5
